$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the missing "Ngày kết thúc" header label in E4, matching the
# formatting already used by the equivalent G4 header cell.
$ws.Range("E4").Value = "Ngày kết thúc"
$ws.Range("G4").Copy()
$ws.Range("E4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the cursor on G7, matching where the author left off.
[void]$ws.Range("G7").Select()
